# Bug fix in Eduati data files (SW837_noCTRL_meas.xlsx)
#
# Sheet1 had 43 stray leftover rows (45:87) below its real A1:N44 data
# table -- each holding nothing but a left-over index value in column A.
# Remove them so the sheet's used range shrinks back down to A1:N44,
# matching Sheet2 / Sheet3.
#
# The workbook is also re-pointed at Sheet1 (it had been left on Sheet3),
# with the selection left sitting further down the (now shorter) sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Drop the stray trailing rows on Sheet1 (rows 45-87 -> dimension back to N44)
$ws1.Rows("45:87").Delete()

# Make Sheet1 the active / tab-selected sheet again (was Sheet3), and leave
# the selection on C56 as recorded by the saved view state.
$ws1.Activate()
$ws1.Range("C56").Select()
